$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the numeric-looking mobile/taluka/village columns
# so they land as shared strings, matching columns A/D/E in existing rows.
$ws.Range("A6:A9").NumberFormat = "@"
$ws.Range("D6:E9").NumberFormat = "@"

# Row 6: ભનુભાઈ જીવરાજભાઈ સોજીત્રા upload
$ws.Range("A6").Value = "9429558759"
$ws.Range("B6").Value = "ભનુભાઈ જીવરાજભાઈ સોજીત્રા"
$ws.Range("C6").Value = "28/05/2025"
$ws.Range("D6").Value = "1"
$ws.Range("E6").Value = "221"

# Row 7: સ્વ. ભનુભાઇ જીવરાજભાઈ સોજિત્રા upload
$ws.Range("A7").Value = "9510851351"
$ws.Range("B7").Value = "સ્વ. ભનુભાઇ જીવરાજભાઈ સોજિત્રા "
$ws.Range("C7").Value = "તા.૨૮-૦૫–૨૦૨૫ વાર બુધવાર "
$ws.Range("D7").Value = "1"
$ws.Range("E7").Value = "221"

# Row 8: duplicate of row 7
$ws.Range("A8").Value = "9510851351"
$ws.Range("B8").Value = "સ્વ. ભનુભાઇ જીવરાજભાઈ સોજિત્રા "
$ws.Range("C8").Value = "તા.૨૮-૦૫–૨૦૨૫ વાર બુધવાર "
$ws.Range("D8").Value = "1"
$ws.Range("E8").Value = "221"

# Row 9: 9429558756 upload
$ws.Range("A9").Value = "9429558756"
$ws.Range("B9").Value = "સ્વ. ભનુભાઈ જીવરાજભાઈ સોજીત્રા|"
$ws.Range("C9").Value = "28/05/2025 | મંગળવાર"
$ws.Range("D9").Value = "1"
$ws.Range("E9").Value = "221"
